# data updated till 15 Dec 11PM
# Column U corresponds to 15-Dec-2020 (G=1-Dec ... U=15-Dec).
# For each retailer row below, record the order received on 15-Dec by
# writing into column U. E (=F/1.04) and F (=SUM(G:AK)) are formulas and
# will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U6").Value  = 2080
$ws.Range("U16").Value = 3120
$ws.Range("U23").Value = 1040
$ws.Range("U30").Value = 1040
$ws.Range("U50").Value = 5200
$ws.Range("U68").Value = 3120
$ws.Range("U78").Value = 5200
$ws.Range("U80").Value = 5200
$ws.Range("U85").Value = 2080
$ws.Range("U90").Value = 2080

# Row 62 is highlighted (matches the existing "latest entry" accent fill
# used elsewhere in the sheet, e.g. T35 / S43) rather than the plain style.
$ws.Range("U62").Value = 2080
$ws.Range("U62").Interior.Color = 8698081

# Reflect where the user was last working when the file was saved.
[void]$ws.Range("U65").Select()
